$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G3").Value = "Yes"

$ws.Range("F7").Value = "Regression"
$ws.Range("G7").Value = "Yes"

$ws.Range("F8").Value = "Contact Us page"
$ws.Range("F9").Value = "Contact Us page"
$ws.Range("F10").Value = "Contact Us page"
$ws.Range("F11").Value = "Contact Us page"
$ws.Range("F12").Value = "Contact Us page"
$ws.Range("F13").Value = "Contact Us page"
$ws.Range("F14").Value = "Contact Us page"
$ws.Range("F15").Value = "Contact Us page"

$ws.Range("F7").Select()
